$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Move the two existing labelled cells (F14 "Mean", F15 "StDev") over to
# column G *before* we drop new SUM() formulas into F14/F15, so the shared
# strings survive (Cut preserves the underlying shared-string cell intact).
# ---------------------------------------------------------------------------
$ws.Range("F14").Cut($ws.Range("G14"))
$ws.Range("F15").Cut($ws.Range("G15"))

# ---------------------------------------------------------------------------
# Highlight fills on D4 (orange) and D8 (green)
# ---------------------------------------------------------------------------
$ws.Range("D4").Interior.Color = 49407
$ws.Range("D8").Interior.Color = 5296274

# ---------------------------------------------------------------------------
# New rolling-5-year SUM column in F, rows 8-16 (F8 orange, F12 green to
# match the new fills applied elsewhere on the same rows)
# ---------------------------------------------------------------------------
$ws.Range("F8").Formula = "=SUM(B4:B8)"
$ws.Range("F8").Interior.Color = 49407

$ws.Range("F9").Formula = "=SUM(B5:B9)"
$ws.Range("F10").Formula = "=SUM(B6:B10)"
$ws.Range("F11").Formula = "=SUM(B7:B11)"

$ws.Range("F12").Formula = "=SUM(B8:B12)"
$ws.Range("F12").Interior.Color = 5296274

$ws.Range("F13").Formula = "=SUM(B9:B13)"
$ws.Range("F14").Formula = "=SUM(B10:B14)"
$ws.Range("F15").Formula = "=SUM(B11:B15)"
$ws.Range("F16").Formula = "=SUM(B12:B16)"

# ---------------------------------------------------------------------------
# New column L (rows 17-27), mirroring column B
# ---------------------------------------------------------------------------
$ws.Range("L17").Value = 165
$ws.Range("L18").Formula = "=B18"
$ws.Range("L19").Formula = "=B19"
$ws.Range("L20").Formula = "=B20"
$ws.Range("L21").Formula = "=B21"
$ws.Range("L22").Formula = "=B22"
$ws.Range("L23").Formula = "=B23"
$ws.Range("L24").Formula = "=B24"
$ws.Range("L25").Formula = "=B25"
$ws.Range("L26").Formula = "=B26"
$ws.Range("L27").Formula = "=B27"

# Running totals in F22:F24
$ws.Range("F22").Formula = "=SUM(B18:B22)"
$ws.Range("F23").Formula = "=F22+B23"
$ws.Range("F24").Formula = "=F23+B11"

# Notes / totals around rows 17-23 (insertion order matters for the shared
# string table layout: flood/800/150 in that order)
$ws.Range("N23").Value = "flood at his birth?"
$ws.Range("P17").Value = "800 years?"
$ws.Range("P18").Value = "150 after"
$ws.Range("P23").Formula = "=SUM(L17:L22)"

# ---------------------------------------------------------------------------
# New little ratio calculation block, rows 30-32
# ---------------------------------------------------------------------------
$ws.Range("E30").Value = 160
$ws.Range("F30").Value = 50

$ws.Range("E31").Formula = "=E30*F32"
$ws.Range("F31").Value = 1300

$ws.Range("F32").Formula = "=F31/F30"

# ---------------------------------------------------------------------------
# Two new source-link rows at the very bottom
# ---------------------------------------------------------------------------
$ws.Range("B54").Value = "https://courses.byui.edu/BOM%20Timeline/html/timeline.html"
$ws.Range("B55").Value = "http://www.bmaf.org/articles/review_this_land_zarahemla__andersen"

# ---------------------------------------------------------------------------
# Scroll / selection state
# ---------------------------------------------------------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 27
$ws.Range("E32").Select()
